$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "캐글/데이콘 경진대회 Baseline을 잡기 위한 optuna + [xgboost, lightgbm, catboost] 패키지 소개"
$ws.Range("E4").Value = "https://teddylee777.github.io/machine-learning/auto-tuning-models"

$ws.Range("D6").Value = "슈퍼짱짱"

$ws.Range("D14").Value = "모두의연구소"

$ws.Range("D19").Value = "꼬낄콘의 분석일지"

$ws.Range("D20").Value = "ai-creator"

$ws.Range("D23").Value = "Be the only one"

$ws.Range("D26").Value = "ai plus(est soft)"

$ws.Range("D28").Value = "로봇이 아닙니다 "

$ws.Range("D32").Value = "데이터과학 삼학년"

$ws.Range("D39").Value = "deadNstreet"

$ws.Range("D42").Value = "IT_notepad"

$ws.Range("D43").Value = "동신한의 조재성"

$ws.Range("D45").Value = "dive-into-ds"
